$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fkey-composite")

# --- Header text: B1 "main_integer_minmax" -> "main_string_enum" ---
$ws.Range("B1").Value = "main_string_enum"

# Column B used to be autosized for the old, longer header; re-autofit-ish width
# for the new (shorter) header text.
$ws.Columns.Item(2).ColumnWidth = 19

# --- Conditional formatting ---
# A2:A1048576 existing rule gains a composite-key membership check against main!A.
$cfA = $ws.Range("A2:A1048576").FormatConditions.Item(1)
$cfA.Formula1 = "=IF(ISBLANK(A2), FALSE, OR(IF(ISNUMBER(A2), INT(A2) <> A2, TRUE), ISNA(MATCH(A2, 'main'!`$A`$2:`$A`$1048576, 0))))"

# B2:B1048576 gets a brand-new rule mirroring the simple-key sheet's pattern.
$cfB = $ws.Range("B2:B1048576").FormatConditions.Add(2, 0, "=IF(ISBLANK(B2), FALSE, ISNA(MATCH(B2, 'main'!`$F`$2:`$F`$1048576, 0)))")
$cfB.Interior.Color = 13551615

# --- Data validations: dropdown lists sourced from main!A and main!F ---
$valA = $ws.Range("A2:A1048576").Validation
$valA.Add(3, 3, 1, "='main'!`$A`$2:`$A`$1048576")
$valA.ShowError = $false
$valA.ErrorTitle = "Invalid value"
$valA.ErrorMessage = "Value must be in dropdown list"

$valB = $ws.Range("B2:B1048576").Validation
$valB.Add(3, 3, 1, "='main'!`$F`$2:`$F`$1048576")
$valB.ShowError = $false
$valB.ErrorTitle = "Invalid value"
$valB.ErrorMessage = "Value must be in dropdown list"
